$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = 19910303
$ws.Range("B2").Value = "zhaoliu"

# Update row 3
$ws.Range("B3").Value = "lisi"

# Update row 4
$ws.Range("B4").Value = "wangwu"

# Update row 5
$ws.Range("A5").Value = 1957
$ws.Range("B5").Value = "maliu"
$ws.Range("C5").Value = "1@qq.com"

# Column C width (bestFit, width 10.5)
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null

# Update selection to D7
$ws.Range("D7").Select()
